# Generate Report for Handback
# Adds a new row (row 4) for file 4a0a440b-cfa1-4237-bee7-2acb77bde0db.md
# ("in sync with en-US") to the Overview, zh-cn and de-de sheets/tables.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, [string]$text) {
    # Force plain text storage (avoids True/False -> boolean coercion and
    # avoids empty strings being dropped entirely).
    if ($text -eq "") {
        $range.Value = "'"
    } else {
        $range.Value = "'" + $text
    }
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

Set-TextValue $wsOv.Range("A4") "4a0a440b-cfa1-4237-bee7-2acb77bde0db.md"
Set-TextValue $wsOv.Range("B4") "e2e\4a0a440b-cfa1-4237-bee7-2acb77bde0db.md"
Set-TextValue $wsOv.Range("C4") ".md"
Set-TextValue $wsOv.Range("E4") "Handed back: in sync with en-US"
Set-TextValue $wsOv.Range("F4") "Handed back: in sync with en-US"
Set-TextValue $wsOv.Range("G4") "2016-08-15 08:59:49"

# B4 looks like the other file-name hyperlink cells (blue/underlined)
$wsOv.Range("B4").Font.Underline = $true
$wsOv.Range("B4").Font.Color = $wsOv.Range("B2").Font.Color

# G4 carries the same date/time display format as the other date cells
$wsOv.Range("G4").NumberFormat = $wsOv.Range("G2").NumberFormat

$wsOv.Hyperlinks.Add($wsOv.Range("B4"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2884f3bf7456e4f06dbd1f108e72d3d8d1216ba/e2e/4a0a440b-cfa1-4237-bee7-2acb77bde0db.md", `
    [Type]::Missing, [Type]::Missing, "e2e\4a0a440b-cfa1-4237-bee7-2acb77bde0db.md")

$loOv = $wsOv.ListObjects.Item(1)
$loOv.Resize($wsOv.Range("A1:G4"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

Set-TextValue $wsZh.Range("A4") "4a0a440b-cfa1-4237-bee7-2acb77bde0db.md"
Set-TextValue $wsZh.Range("B4") ".md"
Set-TextValue $wsZh.Range("C4") "Handed back: in sync with en-US"
Set-TextValue $wsZh.Range("D4") "e2e"
Set-TextValue $wsZh.Range("E4") "ht"
Set-TextValue $wsZh.Range("F4") "True"
Set-TextValue $wsZh.Range("G4") "4a0a440b-cfa1-4237-bee7-2acb77bde0db.2e884f3bf7456e4f06dbd1f108e72d3d8d1216ba.zh-cn.xlf"
Set-TextValue $wsZh.Range("H4") "2016-08-15 08:59:43"
Set-TextValue $wsZh.Range("I4") "4a0a440b-cfa1-4237-bee7-2acb77bde0db.md"
Set-TextValue $wsZh.Range("J4") "4a0a440b-cfa1-4237-bee7-2acb77bde0db.2e884f3bf7456e4f06dbd1f108e72d3d8d1216ba.zh-cn.xlf"
Set-TextValue $wsZh.Range("K4") "2016-08-15 09:00:25"
Set-TextValue $wsZh.Range("L4") ""
Set-TextValue $wsZh.Range("M4") "True"
Set-TextValue $wsZh.Range("N4") ""
Set-TextValue $wsZh.Range("O4") "False"
Set-TextValue $wsZh.Range("P4") ""

$wsZh.Range("A4").Font.Underline = $true
$wsZh.Range("A4").Font.Color = $wsZh.Range("A2").Font.Color
$wsZh.Range("I4").Font.Underline = $true
$wsZh.Range("I4").Font.Color = $wsZh.Range("I2").Font.Color

$wsZh.Range("H4").NumberFormat = $wsZh.Range("H2").NumberFormat
$wsZh.Range("K4").NumberFormat = $wsZh.Range("K2").NumberFormat

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2884f3bf7456e4f06dbd1f108e72d3d8d1216ba/e2e/4a0a440b-cfa1-4237-bee7-2acb77bde0db.md", `
    [Type]::Missing, [Type]::Missing, "4a0a440b-cfa1-4237-bee7-2acb77bde0db.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I4"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e2884f3bf7456e4f06dbd1f108e72d3d8d1216ba/e2e/4a0a440b-cfa1-4237-bee7-2acb77bde0db.md", `
    [Type]::Missing, [Type]::Missing, "4a0a440b-cfa1-4237-bee7-2acb77bde0db.md")

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P4"))

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

Set-TextValue $wsDe.Range("A4") "4a0a440b-cfa1-4237-bee7-2acb77bde0db.md"
Set-TextValue $wsDe.Range("B4") ".md"
Set-TextValue $wsDe.Range("C4") "Handed back: in sync with en-US"
Set-TextValue $wsDe.Range("D4") "e2e"
Set-TextValue $wsDe.Range("E4") "ht"
Set-TextValue $wsDe.Range("F4") "True"
Set-TextValue $wsDe.Range("G4") "4a0a440b-cfa1-4237-bee7-2acb77bde0db.2e884f3bf7456e4f06dbd1f108e72d3d8d1216ba.de-de.xlf"
Set-TextValue $wsDe.Range("H4") "2016-08-15 08:59:49"
Set-TextValue $wsDe.Range("I4") "4a0a440b-cfa1-4237-bee7-2acb77bde0db.md"
Set-TextValue $wsDe.Range("J4") "4a0a440b-cfa1-4237-bee7-2acb77bde0db.2e884f3bf7456e4f06dbd1f108e72d3d8d1216ba.de-de.xlf"
Set-TextValue $wsDe.Range("K4") "2016-08-15 09:00:32"
Set-TextValue $wsDe.Range("L4") ""
Set-TextValue $wsDe.Range("M4") "True"
Set-TextValue $wsDe.Range("N4") ""
Set-TextValue $wsDe.Range("O4") "False"
Set-TextValue $wsDe.Range("P4") ""

$wsDe.Range("A4").Font.Underline = $true
$wsDe.Range("A4").Font.Color = $wsDe.Range("A2").Font.Color
$wsDe.Range("I4").Font.Underline = $true
$wsDe.Range("I4").Font.Color = $wsDe.Range("I2").Font.Color

$wsDe.Range("H4").NumberFormat = $wsDe.Range("H2").NumberFormat
$wsDe.Range("K4").NumberFormat = $wsDe.Range("K2").NumberFormat

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2884f3bf7456e4f06dbd1f108e72d3d8d1216ba/e2e/4a0a440b-cfa1-4237-bee7-2acb77bde0db.md", `
    [Type]::Missing, [Type]::Missing, "4a0a440b-cfa1-4237-bee7-2acb77bde0db.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I4"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/e2884f3bf7456e4f06dbd1f108e72d3d8d1216ba/e2e/4a0a440b-cfa1-4237-bee7-2acb77bde0db.md", `
    [Type]::Missing, [Type]::Missing, "4a0a440b-cfa1-4237-bee7-2acb77bde0db.md")

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P4"))
